$wb = $excel.ActiveWorkbook

# The existing sheet (currently "Sheet1", content "this is in first commit")
# ends up becoming the final "Sheet3" (last tab, still the active/selected
# one) with its single cell updated to the "third commit" string. Two
# brand-new sheets are inserted after it and then the original sheet is
# moved to the very end, so the final tab order is Sheet1(new),
# Sheet2(new), Sheet3(old).
$old = $wb.Worksheets.Item(1)
$tmp1 = $wb.Worksheets.Add($null, $old)
$tmp2 = $wb.Worksheets.Add($null, $tmp1)

# Move the original sheet (still at position 1) to the end.
$old.Move($null, $tmp2) | Out-Null

# Re-resolve fresh references by position now that the move has happened -
# sheet handles obtained before a Move can resolve to the wrong sheet
# afterwards, so always re-fetch via Item() post-move.
$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet3 = $wb.Worksheets.Item(3)

# Rename avoiding transient collisions: $sheet3 already holds the name
# "Sheet1" (the default auto-name) at this point, so give it a throwaway
# name first before the other two claim "Sheet1"/"Sheet2".
$sheet3.Name = "SheetTmp"
$sheet1.Name = "Sheet1"
$sheet2.Name = "Sheet2"
$sheet3.Name = "Sheet3"

# New Sheet1 content.
$sheet1.Range("A1").Value = "this is in first commit"
$sheet1.Range("A2").Value = "this is my second commit"

# New Sheet2 content.
$sheet2.Range("A1").Value = "this is also part of second commit"
$sheet2.Range("A2").Value = "added this in 3rd commit"

# Sheet3 (original sheet) gets its single cell replaced.
$sheet3.Range("A1").Value = "this is part of third commit"

# Match each sheet's saved selection / active cell and make Sheet3 the
# active tab (matches activeTab="2" / tabSelected on the original sheet).
$sheet1.Range("A3").Select() | Out-Null
$sheet2.Range("A3").Select() | Out-Null
$sheet3.Activate() | Out-Null
$sheet3.Range("A2").Select() | Out-Null
